$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = '310.49' }
    @{ Row = 2; Col = 5; Value = '0.23%' }
    @{ Row = 3; Col = 4; Value = '39.66' }
    @{ Row = 3; Col = 5; Value = '-3.51%' }
    @{ Row = 4; Col = 4; Value = '5.073' }
    @{ Row = 4; Col = 5; Value = '-3.50%' }
    @{ Row = 5; Col = 4; Value = '0.07558' }
    @{ Row = 5; Col = 5; Value = '-1.32%' }
    @{ Row = 6; Col = 4; Value = '4.302' }
    @{ Row = 6; Col = 5; Value = '-0.68%' }
    @{ Row = 7; Col = 4; Value = '1.686' }
    @{ Row = 7; Col = 5; Value = '3.88%' }
    @{ Row = 8; Col = 4; Value = '0.9301' }
    @{ Row = 8; Col = 5; Value = '1.23%' }
    @{ Row = 9; Col = 5; Value = '-2.06%' }
    @{ Row = 10; Col = 4; Value = '0.1236' }
    @{ Row = 10; Col = 5; Value = '-0.68%' }
    @{ Row = 11; Col = 4; Value = '0.1789' }
    @{ Row = 11; Col = 5; Value = '-2.77%' }
    @{ Row = 12; Col = 4; Value = '0.09044' }
    @{ Row = 12; Col = 5; Value = '-0.94%' }
    @{ Row = 13; Col = 4; Value = '0.04171' }
    @{ Row = 13; Col = 5; Value = '-4.34%' }
    @{ Row = 14; Col = 4; Value = '0.1053' }
    @{ Row = 14; Col = 5; Value = '0.22%' }
    @{ Row = 15; Col = 4; Value = '0.001288' }
    @{ Row = 15; Col = 5; Value = '2.06%' }
    @{ Row = 16; Col = 4; Value = '0.005870' }
    @{ Row = 16; Col = 5; Value = '1.49%' }
    @{ Row = 18; Col = 4; Value = '3.350' }
    @{ Row = 18; Col = 5; Value = '-0.15%' }
    @{ Row = 19; Col = 4; Value = '0.3353' }
    @{ Row = 19; Col = 5; Value = '0.51%' }
    @{ Row = 20; Col = 4; Value = '7.732' }
    @{ Row = 20; Col = 5; Value = '7.54%' }
    @{ Row = 21; Col = 4; Value = '0.1354' }
    @{ Row = 21; Col = 5; Value = '-2.16%' }
    @{ Row = 22; Col = 4; Value = '0.2892' }
    @{ Row = 22; Col = 5; Value = '-1.09%' }
    @{ Row = 23; Col = 4; Value = '0.04020' }
    @{ Row = 23; Col = 5; Value = '-1.29%' }
    @{ Row = 24; Col = 5; Value = '0.41%' }
    @{ Row = 25; Col = 4; Value = '0.004043' }
    @{ Row = 25; Col = 5; Value = '-3.00%' }
    @{ Row = 26; Col = 5; Value = '0.13%' }
    @{ Row = 38; Col = 4; Value = '0.02419' }
    @{ Row = 38; Col = 5; Value = '-1.18%' }
    @{ Row = 39; Col = 4; Value = '0.05118' }
    @{ Row = 39; Col = 5; Value = '-3.37%' }
    @{ Row = 40; Col = 4; Value = '0.007732' }
    @{ Row = 40; Col = 5; Value = '-1.46%' }
    @{ Row = 41; Col = 4; Value = '0.1295' }
    @{ Row = 41; Col = 5; Value = '-1.48%' }
    @{ Row = 42; Col = 4; Value = '0.007689' }
    @{ Row = 42; Col = 5; Value = '12.82%' }
    @{ Row = 43; Col = 5; Value = '14.32%' }
    @{ Row = 44; Col = 4; Value = '0.008018' }
    @{ Row = 44; Col = 5; Value = '-3.88%' }
    @{ Row = 45; Col = 4; Value = '0.3103' }
    @{ Row = 45; Col = 5; Value = '-7.01%' }
    @{ Row = 46; Col = 4; Value = '0.00006617' }
    @{ Row = 46; Col = 5; Value = '-4.15%' }
    @{ Row = 47; Col = 4; Value = '0.00000000752' }
    @{ Row = 47; Col = 5; Value = '0.02%' }
    @{ Row = 48; Col = 4; Value = '0.2709' }
    @{ Row = 48; Col = 5; Value = '31.84%' }
    @{ Row = 49; Col = 4; Value = '0.004209' }
    @{ Row = 49; Col = 5; Value = '2.67%' }
    @{ Row = 50; Col = 4; Value = '0.00002105' }
    @{ Row = 50; Col = 5; Value = '0.02%' }
    @{ Row = 51; Col = 4; Value = '0.0002004' }
    @{ Row = 51; Col = 5; Value = '0.02%' }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
